$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.989.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.12%  "

$ws.Range("D3").Value = "'2.464.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.88%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'516.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.85%  "

$ws.Range("D6").Value = "'131.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.47%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "'0.555"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.14%  "

$ws.Range("D9").Value = "'0.0988"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.63%  "

$ws.Range("E10").Value = "  -1.77%  "

$ws.Range("D11").Value = "'5.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.72%  "

$ws.Range("D12").Value = "'0.340"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.42%  "

$ws.Range("D13").Value = "'2.901.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.95%  "

$ws.Range("D14").Value = "'57.893.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.11%  "

$ws.Range("D15").Value = "'21.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.07%  "

$ws.Range("E16").Value = "  -2.74%  "

$ws.Range("D17").Value = "'2.465.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.50%  "

$ws.Range("D18").Value = "'10.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.90%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'319.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.01%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.22%  "

$ws.Range("D21").Value = "'0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").Value = "'5.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.78%  "

$ws.Range("D23").Value = "'64.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.14%  "

$ws.Range("D24").Value = "'0.407"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.96%  "

$ws.Range("E25").Value = "  -0.35%  "

$ws.Range("D26").Value = "'0.161"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.43%  "

$ws.Range("D27").Value = "'7.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.79%  "

$ws.Range("D28").Value = "'0.0₃0744"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.33%  "

$ws.Range("D29").Value = "'6.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.40%  "

$ws.Range("D30").Value = "'1.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.01%  "

$ws.Range("D31").Value = "'165.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("E32").Value = "  -4.22%  "

$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("D34").Value = "'0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("D35").Value = "'18.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.28%  "

$ws.Range("E36").Value = "  -9.64%  "

$ws.Range("E37").Value = "  -4.46%  "

$ws.Range("D38").Value = "'1.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.88%  "

$ws.Range("D39").Value = "'0.786"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.03%  "

$ws.Range("D40").Value = "'3.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.93%  "

$ws.Range("D41").Value = "'273.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.65%  "

$ws.Range("D42").Value = "'4.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.66%  "

$ws.Range("D43").Value = "'0.588"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.64%  "

$ws.Range("D44").Value = "'126.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.89%  "

$ws.Range("D45").Value = "'0.0906"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.87%  "

$ws.Range("D46").Value = "'0.0489"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.21%  "

$ws.Range("E47").Value = "  -3.54%  "

$ws.Range("D48").Value = "'16.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.46%  "

$ws.Range("D49").Value = "'1.729.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.94%  "

$ws.Range("E50").Value = "  -1.77%  "

$ws.Range("E51").Value = "  -2.53%  "

Write-Output "applied changes"